# Applies the diff: swap the contents/styles of columns C and D (the
# " Charu Python Average Power(W) after taking abs" column now comes
# before the " Charu Python Average Power(W)" column), widen/resize
# columns C & D, add a new (empty) column G formatted like column C,
# and update the view (zoom + selected cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap columns C and D -------------------------------------------------
# Cutting column D and inserting the cut cells in front of column C moves
# D's whole content+formatting into C, and shifts the old C into D -
# exactly mirroring the header/value/style swap seen in the diff.
$ws.Columns("D").Cut() | Out-Null
$ws.Columns("C").Insert() | Out-Null

# --- Column widths ---------------------------------------------------------
$ws.Columns("C").ColumnWidth = 50.16666666
$ws.Columns("D").ColumnWidth = 39.66666666

# --- New column G ------------------------------------------------------
# Column G is a brand-new, empty column that only carries the same
# per-row formatting as the (new) column C.
$ws.Range("C1:C9").Copy() | Out-Null
$ws.Range("G1:G9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Sheet view: zoom + selection -------------------------------------
$ws.Range("C12").Select() | Out-Null
$excel.ActiveWindow.Zoom = 85
